# Scheduled market-price refresh: update currentAveragePrice / Leve profit
# columns (H:N) for the rows whose underlying item prices changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 55724
$ws.Range("I33").Value = 31405.719
$ws.Range("J33").Value = 250270.25
$ws.Range("K33").Value = 31405.719
$ws.Range("L33").Value = 250270.25
$ws.Range("M33").Value = -31176.719
$ws.Range("N33").Value = -250728.25

$ws.Range("H38").Value = 541.3125
$ws.Range("I38").Value = 225.6923
$ws.Range("J38").Value = 1909
$ws.Range("K38").Value = 677.0769
$ws.Range("L38").Value = 5727
$ws.Range("M38").Value = -305.0769
$ws.Range("N38").Value = -6471

$ws.Range("H62").Value = 2308.25
$ws.Range("I62").Value = 1950
$ws.Range("J62").Value = 2666.5
$ws.Range("K62").Value = 1950
$ws.Range("L62").Value = 2666.5
$ws.Range("M62").Value = -1326
$ws.Range("N62").Value = -3914.5

$ws.Range("H65").Value = 2308.25
$ws.Range("I65").Value = 1950
$ws.Range("J65").Value = 2666.5
$ws.Range("K65").Value = 9750
$ws.Range("L65").Value = 13332.5
$ws.Range("M65").Value = -6630
$ws.Range("N65").Value = -19572.5

$ws.Range("H70").Value = 2220.2
$ws.Range("I70").Value = 2340.4
$ws.Range("J70").Value = 2100
$ws.Range("K70").Value = 7021.200000000001
$ws.Range("L70").Value = 6300
$ws.Range("M70").Value = -6751.200000000001
$ws.Range("N70").Value = -6840

$ws.Range("H73").Value = 2220.2
$ws.Range("I73").Value = 2340.4
$ws.Range("J73").Value = 2100
$ws.Range("K73").Value = 7021.200000000001
$ws.Range("L73").Value = 6300
$ws.Range("M73").Value = -6085.200000000001
$ws.Range("N73").Value = -8172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1601
$ws.Range("I45").Value = 1639.7
$ws.Range("J45").Value = 1214
$ws.Range("K45").Value = 1639.7
$ws.Range("L45").Value = 1214
$ws.Range("M45").Value = -1262.7
$ws.Range("N45").Value = -1968

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1926.8846
$ws.Range("I86").Value = 2154.95
$ws.Range("J86").Value = 1166.6666
$ws.Range("K86").Value = 2154.95
$ws.Range("L86").Value = 1166.6666
$ws.Range("M86").Value = -1031.95
$ws.Range("N86").Value = -3412.6666

$ws.Range("H89").Value = 1926.8846
$ws.Range("I89").Value = 2154.95
$ws.Range("J89").Value = 1166.6666
$ws.Range("K89").Value = 10774.75
$ws.Range("L89").Value = 5833.333000000001
$ws.Range("M89").Value = -5158.75
$ws.Range("N89").Value = -17065.333

$ws.Range("H135").Value = 21742.223
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 21742.223
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 21742.223
$ws.Range("N135").Value = -31882.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 24688.857
$ws.Range("I50").Value = 18000
$ws.Range("J50").Value = 25203.385
$ws.Range("K50").Value = 18000
$ws.Range("L50").Value = 25203.385
$ws.Range("M50").Value = -17375
$ws.Range("N50").Value = -26453.385

$ws.Range("H51").Value = 22939.6
$ws.Range("I51").Value = 15000
$ws.Range("J51").Value = 23357.475
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 23357.475
$ws.Range("M51").Value = -14264
$ws.Range("N51").Value = -24829.475

$ws.Range("H59").Value = 42680.285
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 42680.285
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 42680.285
$ws.Range("N59").Value = -44970.285

$ws.Range("H60").Value = 22330.354
$ws.Range("I60").Value = 15000
$ws.Range("J60").Value = 22788.5
$ws.Range("K60").Value = 15000
$ws.Range("L60").Value = 22788.5
$ws.Range("M60").Value = -14489
$ws.Range("N60").Value = -23810.5

$ws.Range("H61").Value = 22939.6
$ws.Range("I61").Value = 15000
$ws.Range("J61").Value = 23357.475
$ws.Range("K61").Value = 15000
$ws.Range("L61").Value = 23357.475
$ws.Range("M61").Value = -14652
$ws.Range("N61").Value = -24053.475

$ws.Range("H68").Value = 41911.668
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 41911.668
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 41911.668
$ws.Range("N68").Value = -43409.668

$ws.Range("H71").Value = 41911.668
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 41911.668
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 125735.004
$ws.Range("N71").Value = -133223.004

$ws.Range("H74").Value = 34436.266
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 34436.266
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 34436.266
$ws.Range("N74").Value = -36184.266

$ws.Range("H77").Value = 34436.266
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 34436.266
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 103308.798
$ws.Range("N77").Value = -112044.798

$ws.Range("H99").Value = 1437
$ws.Range("I99").Value = 1304.5454
$ws.Range("J99").Value = 1645.1428
$ws.Range("K99").Value = 1304.5454
$ws.Range("L99").Value = 1645.1428
$ws.Range("M99").Value = 193.4546
$ws.Range("N99").Value = -4641.1428

$ws.Range("H126").Value = 1437
$ws.Range("I126").Value = 1304.5454
$ws.Range("J126").Value = 1645.1428
$ws.Range("K126").Value = 3913.6362
$ws.Range("L126").Value = 4935.428400000001
$ws.Range("M126").Value = -1443.6362
$ws.Range("N126").Value = -9875.428400000001

$ws.Range("H132").Value = 7152.6177
$ws.Range("I132").Value = 4470.107
$ws.Range("J132").Value = 19671
$ws.Range("K132").Value = 13410.321
$ws.Range("L132").Value = 59013
$ws.Range("M132").Value = -10880.321
$ws.Range("N132").Value = -64073

$ws.Range("M141").ClearContents()
$ws.Range("H141").Value = 50000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 50000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 41.076923
$ws.Range("I8").Value = 41.076923
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 123.230769
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 15.76923099999999

$ws.Range("H100").Value = 3066.7
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3066.7
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 9200.099999999999
$ws.Range("N100").Value = -10822.1

$ws.Range("H103").Value = 3012.8333
$ws.Range("I103").Value = 1875.6923
$ws.Range("J103").Value = 3882.4119
$ws.Range("K103").Value = 5627.0769
$ws.Range("L103").Value = 11647.2357
$ws.Range("M103").Value = -4748.0769
$ws.Range("N103").Value = -13405.2357

$ws.Range("H106").Value = 3500
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 3500
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 10500
$ws.Range("N106").Value = -12392

$ws.Range("H130").Value = 1776.6666
$ws.Range("I130").Value = 415
$ws.Range("J130").Value = 4500
$ws.Range("K130").Value = 1245
$ws.Range("L130").Value = 13500
$ws.Range("M130").Value = 3775
$ws.Range("N130").Value = -23540

$ws.Range("H131").Value = 824.18335
$ws.Range("I131").Value = 351.14285
$ws.Range("J131").Value = 886.6604
$ws.Range("K131").Value = 1053.42855
$ws.Range("L131").Value = 2659.9812
$ws.Range("M131").Value = 3986.57145
$ws.Range("N131").Value = -12739.9812

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2754.8572
$ws.Range("I102").Value = 2456.8
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 2456.8
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -834.8000000000002
$ws.Range("N102").Value = -6744

$ws.Range("H126").Value = 101131.1
$ws.Range("I126").Value = 167561.83
$ws.Range("J126").Value = 1485
$ws.Range("K126").Value = 502685.49
$ws.Range("L126").Value = 4455
$ws.Range("M126").Value = -500215.49
$ws.Range("N126").Value = -9395

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").ClearContents()
$ws.Range("H7").Value = 2151.25
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2151.25
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2151.25
$ws.Range("N7").Value = -2375.25

$ws.Range("H40").Value = 4123
$ws.Range("I40").Value = 4285.923
$ws.Range("J40").Value = 2005
$ws.Range("K40").Value = 4285.923
$ws.Range("L40").Value = 2005
$ws.Range("M40").Value = -4149.923
$ws.Range("N40").Value = -2277

$ws.Range("M126").ClearContents()
$ws.Range("H126").Value = 2151.25
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2151.25
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 6453.75
$ws.Range("N126").Value = -11393.75

$ws.Range("H141").Value = 69203
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 69203
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 69203
$ws.Range("N141").Value = -79563

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 40672.727
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 40672.727
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 40672.727
$ws.Range("N116").Value = -49850.727

$ws.Range("H126").Value = 2480877.2
$ws.Range("I126").Value = 1990580.5
$ws.Range("J126").Value = 3624903.2
$ws.Range("K126").Value = 5971741.5
$ws.Range("L126").Value = 10874709.6
$ws.Range("M126").Value = -5969271.5
$ws.Range("N126").Value = -10879649.6
